$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.56
$ws.Range("I2").Value = 5.75
$ws.Range("Q2").Value = 1.92
$ws.Range("R2").Value = 1.98
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("AI2").Value = 29
$ws.Range("AO2").Value = 8
$ws.Range("AS2").Value = 151
$ws.Range("G3").Value = 1.76
$ws.Range("H3").Value = 3.25
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("AD3").Value = 6.5
$ws.Range("AH3").Value = 9.5
$ws.Range("AK3").Value = 67
$ws.Range("AV3").Value = 101
$ws.Range("G4").Value = 1.71
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("G10").Value = 2.1
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 3.6
$ws.Range("L10").Value = 4.33
$ws.Range("N10").Value = 7.5
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.73
$ws.Range("X10").Value = 9
$ws.Range("AC10").Value = 7.5
$ws.Range("AF10").Value = 67
$ws.Range("AL10").Value = 34
$ws.Range("AO10").Value = 12
$ws.Range("AW10").Value = 5.5
$ws.Range("AY10").Value = 34
$ws.Range("BB10").Value = 301
$ws.Range("G11").Value = 1.9
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 4.33
$ws.Range("J11").Value = 2.75
$ws.Range("L11").Value = 5.5
$ws.Range("U11").Value = 2.5
$ws.Range("V11").Value = 1.5
$ws.Range("X11").Value = 7.5
$ws.Range("Z11").Value = 15
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 17
$ws.Range("AK11").Value = 51
$ws.Range("AN11").Value = 3.6
$ws.Range("AO11").Value = 11
$ws.Range("AX11").Value = 29
$ws.Range("G13").Value = 1.29
$ws.Range("H13").Value = 5
$ws.Range("K13").Value = 2.6
$ws.Range("L13").Value = 8.5
$ws.Range("X13").Value = 6.5
$ws.Range("AE13").Value = 21
$ws.Range("AJ13").Value = 29
$ws.Range("AK13").Value = 126
$ws.Range("AL13").Value = 67
$ws.Range("AN13").Value = 3.25
$ws.Range("AO13").Value = 6
$ws.Range("AP13").Value = 17
